$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 1147.5714
$ws.Range("I39").Value = 588.6
$ws.Range("K39").Value = 1765.8
$ws.Range("M39").Value = -1469.8

$ws.Range("H100").Value = 2154.5715
$ws.Range("I100").Value = 2154.5715
$ws.Range("K100").Value = 2154.5715
$ws.Range("M100").Value = -1613.5715

$ws.Range("H134").Value = 79999.5
$ws.Range("J134").Value = 79999.5
$ws.Range("L134").Value = 79999.5
$ws.Range("N134").Value = -90139.5

$ws.Range("H138").Value = 3822.7693
$ws.Range("J138").Value = 3828.5715
$ws.Range("L138").Value = 11485.7145
$ws.Range("N138").Value = -21765.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3824.25
$ws.Range("I2").Value = 1766.3334
$ws.Range("K2").Value = 1766.3334
$ws.Range("M2").Value = -1653.3334

$ws.Range("H32").Value = 3948.0833
$ws.Range("I32").Value = 3694.7827
$ws.Range("K32").Value = 3694.7827
$ws.Range("M32").Value = -3407.7827

$ws.Range("H45").Value = 3853.0833
$ws.Range("I45").Value = 1449.6666
$ws.Range("J45").Value = 4654.222
$ws.Range("K45").Value = 1449.6666
$ws.Range("L45").Value = 4654.222
$ws.Range("M45").Value = -1072.6666
$ws.Range("N45").Value = -5408.222

$ws.Range("H116").Value = 3824.25
$ws.Range("I116").Value = 1766.3334
$ws.Range("K116").Value = 1766.3334
$ws.Range("M116").Value = 527.6666

$ws.Range("H132").Value = 789.9286
$ws.Range("I132").Value = 789.9286
$ws.Range("K132").Value = 2369.7858
$ws.Range("M132").Value = 160.2142000000003

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3824.25
$ws.Range("I3").Value = 1766.3334
$ws.Range("K3").Value = 1766.3334
$ws.Range("M3").Value = -1652.3334

$ws.Range("H20").Value = 1521.1
$ws.Range("I20").Value = 788.8570999999999
$ws.Range("K20").Value = 788.8570999999999
$ws.Range("M20").Value = -541.8570999999999

$ws.Range("H22").Value = 754.4286
$ws.Range("I22").Value = 939.8
$ws.Range("J22").Value = 291
$ws.Range("K22").Value = 939.8
$ws.Range("L22").Value = 291
$ws.Range("M22").Value = -766.8
$ws.Range("N22").Value = -637

$ws.Range("H81").Value = 17500
$ws.Range("J81").Value = 17500
$ws.Range("L81").Value = 17500
$ws.Range("N81").Value = -19622

$ws.Range("H84").Value = 17500
$ws.Range("J84").Value = 17500
$ws.Range("L84").Value = 52500
$ws.Range("N84").Value = -63108

$ws.Range("H134").Value = 1629.5
$ws.Range("I134").Value = 1629.5
$ws.Range("K134").Value = 4888.5
$ws.Range("M134").Value = -2353.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 1002
$ws.Range("I3").Value = 1002
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1002
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -889
$ws.Range("N3").ClearContents()

$ws.Range("H6").Value = 400.25
$ws.Range("I6").Value = 799.5
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 799.5
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = -686.5
$ws.Range("N6").Value = -227

$ws.Range("H7").Value = 1861.3684
$ws.Range("I7").Value = 1075.6
$ws.Range("K7").Value = 1075.6
$ws.Range("M7").Value = -962.5999999999999

$ws.Range("H10").Value = 1839.5
$ws.Range("J10").Value = 5003
$ws.Range("L10").Value = 5003
$ws.Range("N10").Value = -5281

$ws.Range("H31").Value = 2127.2
$ws.Range("I31").Value = 1045.3334
$ws.Range("J31").Value = 3750
$ws.Range("K31").Value = 1045.3334
$ws.Range("L31").Value = 3750
$ws.Range("M31").Value = -750.3334
$ws.Range("N31").Value = -4340

$ws.Range("H34").Value = 2127.2
$ws.Range("I34").Value = 1045.3334
$ws.Range("J34").Value = 3750
$ws.Range("K34").Value = 1045.3334
$ws.Range("L34").Value = 3750
$ws.Range("M34").Value = -843.3334
$ws.Range("N34").Value = -4154

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 456
$ws.Range("J9").Value = 456
$ws.Range("L9").Value = 1368
$ws.Range("N9").Value = -1816

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 1000
$ws.Range("J4").Value = 1000
$ws.Range("L4").Value = 1000
$ws.Range("N4").Value = -1224

$ws.Range("H62").Value = 52000
$ws.Range("I62").Value = 44000
$ws.Range("K62").Value = 44000
$ws.Range("M62").Value = -43314

$ws.Range("H65").Value = 52000
$ws.Range("I65").Value = 44000
$ws.Range("K65").Value = 132000
$ws.Range("M65").Value = -128568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 2310.6667
$ws.Range("I32").Value = 2310.6667
$ws.Range("K32").Value = 2310.6667
$ws.Range("M32").Value = -1993.6667

$ws.Range("H100").Value = 2366.6667
$ws.Range("I100").Value = 1550
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 1550
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -1009
$ws.Range("N100").Value = -5082

$ws.Range("H122").Value = 5640
$ws.Range("I122").Value = 5640
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 16920
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -14470
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10001.5
$ws.Range("I62").Value = 10000
$ws.Range("J62").Value = 10003
$ws.Range("K62").Value = 10000
$ws.Range("L62").Value = 10003
$ws.Range("M62").Value = -9376
$ws.Range("N62").Value = -11251

$ws.Range("H65").Value = 10001.5
$ws.Range("I65").Value = 10000
$ws.Range("J65").Value = 10003
$ws.Range("K65").Value = 50000
$ws.Range("L65").Value = 50015
$ws.Range("M65").Value = -46880
$ws.Range("N65").Value = -56255

$ws.Range("H81").Value = 3109.8
$ws.Range("I81").Value = 2899.7778
$ws.Range("K81").Value = 5799.5556
$ws.Range("M81").Value = -4738.5556

$ws.Range("H84").Value = 3109.8
$ws.Range("I84").Value = 2899.7778
$ws.Range("K84").Value = 28997.778
$ws.Range("M84").Value = -23693.778

$ws.Range("H112").Value = 30972.5
$ws.Range("J112").Value = 30972.5
$ws.Range("L112").Value = 30972.5
$ws.Range("N112").Value = -33926.5

$ws.Range("H122").Value = 1234.4
$ws.Range("I122").Value = 1234.4
$ws.Range("K122").Value = 3703.2
$ws.Range("M122").Value = -1253.2
